# Adapt column header formatting to respective input file names:
#   "..._old" -> "..._FV2210"
#   "..._new" -> "..._FV2304"
# then (re)build the Excel Table over the data range and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# --- Rename the header row: "_old" suffix -> "_FV2210", "_new" suffix -> "_FV2304" ---
for ($col = 1; $col -le $lastCol; $col++) {
    $headerCell = $ws.Cells.Item(1, $col)
    $headerValue = $headerCell.Value2
    if ($headerValue -ne $null) {
        $renamed = $headerValue -replace "_old$", "_FV2210"
        $renamed = $renamed -replace "_new$", "_FV2304"
        if ($renamed -ne $headerValue) {
            $headerCell.Value2 = $renamed
        }
    }
}

# --- Turn the data range into a native Excel Table (ListObject) with an AutoFilter ---
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# --- Freeze the header row (split after row 1) ---
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Header columns renamed, Table1 created, header row frozen."
